$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.419.54"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "2.984.39"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.48"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "2.983.87"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("E15").Value = "  +3.12%  "
$ws.Range("D16").Value = "3.476.52"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("D18").Value = "61.393.29"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "2.986.40"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "445.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "0.0₃0809"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("E37").Value = "  -0.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("E41").Value = "  +8.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "387.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.269"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0349"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("D47").Value = "2.684.49"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.87%  "
